$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4885.2856
$ws.Range("J29").Value = 4885.2856
$ws.Range("L29").Value = 14655.8568
$ws.Range("N29").Value = -15217.8568

$ws.Range("H38").Value = 3798.5557
$ws.Range("I38").Value = 598.2857
$ws.Range("J38").Value = 14999.5
$ws.Range("K38").Value = 1794.8571
$ws.Range("L38").Value = 44998.5
$ws.Range("M38").Value = -1422.8571
$ws.Range("N38").Value = -45742.5

$ws.Range("H40").Value = 1849.875
$ws.Range("I40").Value = 1833.1666
$ws.Range("J40").Value = 1900
$ws.Range("K40").Value = 1833.1666
$ws.Range("L40").Value = 1900
$ws.Range("M40").Value = -1658.1666
$ws.Range("N40").Value = -2250

$ws.Range("H58").Value = 3026.6667
$ws.Range("I58").Value = 53.333332
$ws.Range("K58").Value = 159.999996
$ws.Range("M58").Value = -9.99999600000001

$ws.Range("H64").Value = 4942.857
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 4920
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 4920
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -5416

$ws.Range("H67").Value = 4942.857
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 4920
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 4920
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -6636

$ws.Range("H131").Value = 7526.8887
$ws.Range("I131").Value = 2580.6667
$ws.Range("K131").Value = 7742.000100000001
$ws.Range("M131").Value = -2702.000100000001

$ws.Range("H132").Value = 1929.1786
$ws.Range("I132").Value = 1851.2273
$ws.Range("K132").Value = 5553.6819
$ws.Range("M132").Value = -3023.6819

$ws.Range("H133").Value = 206311.6
$ws.Range("J133").Value = 206311.6
$ws.Range("L133").Value = 206311.6
$ws.Range("N133").Value = -216431.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21993.834
$ws.Range("I82").Value = 6392.8
$ws.Range("K82").Value = 6392.8
$ws.Range("M82").Value = -6009.8

$ws.Range("H85").Value = 21993.834
$ws.Range("I85").Value = 6392.8
$ws.Range("K85").Value = 6392.8
$ws.Range("M85").Value = -5066.8

$ws.Range("H97").Value = 13189.667
$ws.Range("I97").Value = 13189.667
$ws.Range("K97").Value = 13189.667
$ws.Range("M97").Value = -12198.667

$ws.Range("H107").Value = 2036.909
$ws.Range("I107").Value = 1434
$ws.Range("K107").Value = 1434
$ws.Range("M107").Value = 486

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4062.9285
$ws.Range("J31").Value = 5015.9
$ws.Range("L31").Value = 5015.9
$ws.Range("N31").Value = -5605.9

$ws.Range("H34").Value = 4062.9285
$ws.Range("J34").Value = 5015.9
$ws.Range("L34").Value = 5015.9
$ws.Range("N34").Value = -5419.9

$ws.Range("H42").Value = 4519
$ws.Range("I42").Value = 38
$ws.Range("J42").Value = 9000
$ws.Range("K42").Value = 38
$ws.Range("L42").Value = 9000
$ws.Range("M42").Value = 555
$ws.Range("N42").Value = -10186

$ws.Range("H44").Value = 16000
$ws.Range("J44").Value = 16000
$ws.Range("L44").Value = 16000
$ws.Range("N44").Value = -16884

$ws.Range("H55").Value = 10081
$ws.Range("J55").Value = 10081
$ws.Range("L55").Value = 10081
$ws.Range("N55").Value = -10711

$ws.Range("H62").Value = 44311.1
$ws.Range("J62").Value = 83896.8
$ws.Range("L62").Value = 83896.8
$ws.Range("N62").Value = -85144.8

$ws.Range("H65").Value = 44311.1
$ws.Range("J65").Value = 83896.8
$ws.Range("L65").Value = 419484
$ws.Range("N65").Value = -425724

$ws.Range("H132").Value = 4084.9473
$ws.Range("I132").Value = 3180.0833
$ws.Range("K132").Value = 9540.249899999999
$ws.Range("M132").Value = -7010.249899999999

$ws.Range("H134").Value = 4922.9414
$ws.Range("I134").Value = 3606.7273
$ws.Range("K134").Value = 10820.1819
$ws.Range("M134").Value = -8285.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5882577
$ws.Range("I7").Value = 9091004
$ws.Range("K7").Value = 27273012
$ws.Range("M7").Value = -27272900

$ws.Range("H80").Value = 5864.222
$ws.Range("I80").Value = 5796.3335
$ws.Range("K80").Value = 17389.0005
$ws.Range("M80").Value = -16453.0005

$ws.Range("H83").Value = 5864.222
$ws.Range("I83").Value = 5796.3335
$ws.Range("K83").Value = 52167.0015
$ws.Range("M83").Value = -47487.0015

$ws.Range("H88").Value = 5407
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15856

$ws.Range("H91").Value = 5407
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17964

$ws.Range("H129").Value = 3872.8
$ws.Range("J129").Value = 4000
$ws.Range("L129").Value = 12000
$ws.Range("N129").Value = -22000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 295.93332
$ws.Range("J2").Value = 412.375
$ws.Range("L2").Value = 412.375
$ws.Range("N2").Value = -638.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H132").Value = 5374
$ws.Range("J132").Value = 6249.5293
$ws.Range("L132").Value = 18748.5879
$ws.Range("N132").Value = -23808.5879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7105.5
$ws.Range("J62").Value = 7506.3125
$ws.Range("L62").Value = 7506.3125
$ws.Range("N62").Value = -8754.3125

$ws.Range("H65").Value = 7105.5
$ws.Range("J65").Value = 7506.3125
$ws.Range("L65").Value = 37531.5625
$ws.Range("N65").Value = -43771.5625

$ws.Range("H81").Value = 1820.3334
$ws.Range("I81").Value = 1820.3334
$ws.Range("K81").Value = 3640.6668
$ws.Range("M81").Value = -2579.6668

$ws.Range("H84").Value = 1820.3334
$ws.Range("I84").Value = 1820.3334
$ws.Range("K84").Value = 18203.334
$ws.Range("M84").Value = -12899.334

$ws.Range("H126").Value = 171250.5
$ws.Range("I126").Value = 202500.6
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 607501.8
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -605031.8
$ws.Range("N126").Value = -49940

$ws.Range("H132").Value = 1658.5
$ws.Range("I132").Value = 1145.4667
$ws.Range("K132").Value = 3436.4001
$ws.Range("M132").Value = -906.4000999999998

$ws.Range("H136").Value = 50690.24
$ws.Range("I136").Value = 2677.6875
$ws.Range("K136").Value = 8033.0625
$ws.Range("M136").Value = -5483.0625
